$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text (they use "." as thousands sep, %, etc.)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "60.820.79"
$ws.Range("E2").Value = "  -3.29%  "
$ws.Range("D3").Value = "2.913.41"
$ws.Range("E3").Value = "  -3.98%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "583.80"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").Value = "144.85"
$ws.Range("E6").Value = "  -5.68%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").Value = "2.911.78"
$ws.Range("E9").Value = "  -3.84%  "
$ws.Range("D10").Value = "6.83"
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("E11").Value = "  -4.64%  "
$ws.Range("E12").Value = "  -4.02%  "
$ws.Range("E13").Value = "  -3.86%  "
$ws.Range("D14").Value = "33.49"
$ws.Range("E14").Value = "  -5.99%  "
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "3.394.96"
$ws.Range("E16").Value = "  -4.03%  "
$ws.Range("D17").Value = "60.750.37"
$ws.Range("E17").Value = "  -3.33%  "
$ws.Range("E18").Value = "  -5.45%  "
$ws.Range("D19").Value = "2.906.47"
$ws.Range("E19").Value = "  -4.25%  "
$ws.Range("D20").Value = "431.80"
$ws.Range("E20").Value = "  -4.66%  "
$ws.Range("D21").Value = "13.61"
$ws.Range("E21").Value = "  -4.69%  "
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("E23").Value = "  -4.64%  "
$ws.Range("D24").Value = "80.32"
$ws.Range("E24").Value = "  -3.48%  "
$ws.Range("D25").Value = "10.84"
$ws.Range("E25").Value = "  -4.08%  "
$ws.Range("E26").Value = "  -4.83%  "
$ws.Range("D27").Value = "11.86"
$ws.Range("E27").Value = "  -4.20%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "7.17"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("D31").Value = "2.61"
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("D33").Value = "26.50"
$ws.Range("E33").Value = "  -3.86%  "
$ws.Range("E34").Value = "  -3.81%  "
$ws.Range("D35").Value = "0.0₃0865"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("D37").Value = "5.65"
$ws.Range("E37").Value = "  -4.97%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.128"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "3.02"
$ws.Range("E39").Value = "  -5.51%  "
$ws.Range("D40").Value = "49.74"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("E41").Value = "  -5.23%  "
$ws.Range("D42").Value = "8.65"
$ws.Range("E42").Value = "  -4.91%  "
$ws.Range("D43").Value = "0.293"
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("D44").Value = "41.16"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").Value = "375.45"
$ws.Range("E45").Value = "  -5.37%  "
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("D47").Value = "2.667.55"
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D48").Value = "132.71"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D50").Value = "24.28"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("E51").Value = "  -1.92%  "
